# Weekly update: insert a new price-report row (new market day) at row 4,
# pushing the existing rows 4-77 down to 5-78.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 4 (shifts rows 4:77 down to 5:78,
# dimension grows to A1:R78 automatically).
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with the new weekly record.
$ws.Range("A4").Value2 = 2
$ws.Range("B4").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C4").Value2 = "Coquimbo"
$ws.Range("D4").Value2 = 44860
$ws.Range("E4").Value2 = 4
$ws.Range("F4").Value2 = 100112026
$ws.Range("G4").Value2 = "Haba"
$ws.Range("H4").Value2 = "Sin especificar"
$ws.Range("I4").Value2 = "Primera"
$ws.Range("J4").Value2 = 700
$ws.Range("K4").Value2 = 4000
$ws.Range("L4").Value2 = 5000
$ws.Range("M4").Value2 = 4500
$ws.Range("N4").Value2 = "$/saco 25 kilos"
$ws.Range("O4").Value2 = "Provincia de Limarí"
$ws.Range("P4").Value2 = 180
$ws.Range("Q4").Value2 = 25
$ws.Range("R4").Value2 = "Hortaliza"
